$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.778.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.865.94'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.23'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7006'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07770'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3083'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.83'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07843'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.186'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.867.90'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.82'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6962'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.652'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.771.20'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008393'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.116.82'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '243.88'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.83'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.653'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.973'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.40'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.43'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.544'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.286'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.218'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.200'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05099'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7904'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.925'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.164'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.338.31'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +9.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01885'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.755'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9668'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.048'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +12.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.95'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +6.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.797'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.015.38'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.53'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.92%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5200'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.55%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.797'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.037'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.04%  '